# Daily attendance processing - 2026-01-08 10:07:15
# Normalizes the "Recorded By" (column G) entries so that the "System"
# token is moved from the front of the comma-separated list to the end,
# and is written with proper capitalization ("System") once relocated.
#
# Examples:
#   "System, dnasr281@gmail.com"                 -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"                -> "backup@backdoor.com, System"
#   "System, backup@backdoor.com, system"        -> "system, backup@backdoor.com, System"
#
# Entries that do not start with a "System"/"system" token (or that
# consist of a single token only) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$colIndex = 7  # column G = "Recorded By"
$updated = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $rawParts = $text.Split(",")
    if ($rawParts.Length -lt 2) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts[0].ToLower() -ne "system") {
        continue
    }

    if ($parts.Length -eq 2) {
        $newVal = $parts[1] + ", System"
    } elseif ($parts.Length -eq 3) {
        $newVal = $parts[2] + ", " + $parts[1] + ", System"
    } else {
        continue
    }

    if (-not $newVal.Equals($text)) {
        $cell.Value = $newVal
        $updated++
    }
}

Write-Host ("Updated " + $updated + " 'Recorded By' cell(s).")
